# Disease Maps logo deck — apply the authored edit:
#   * Reorder slides: the grouped "disease maps" logo slide moves to the
#     front, a brand-new variant logo slide is inserted as slide 2, and the
#     old standalone-textbox logo slide becomes an (emptied) slide 3.
#   * The remaining two slides simply shift down one position.

function HexColor($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Move the grouped logo slide (currently 2nd) to the front.
# ---------------------------------------------------------------------
$p.Slides.Item(2).MoveTo(1)

# ---------------------------------------------------------------------
# 2) The old front slide (standalone 4-textbox logo) is now 2nd; its
#    content is removed, leaving an empty slide (it keeps its identity /
#    position further down the deck once the new slide is inserted).
# ---------------------------------------------------------------------
$emptied = $p.Slides.Item(2)
for ($i = $emptied.Shapes.Count; $i -ge 1; $i--) {
    $emptied.Shapes.Item($i).Delete()
}

# ---------------------------------------------------------------------
# 3) Insert a brand-new slide at position 2, using the same blank custom
#    layout as the other logo slides.
# ---------------------------------------------------------------------
$layout = $p.SlideMaster.CustomLayouts.Item(5)
$newSlide = $p.Slides.AddSlide(2, $layout)

# ---------------------------------------------------------------------
# 4) Populate the new slide with the repositioned logo word-marks plus a
#    small rounded-rectangle accent, matching the authored layout.
# ---------------------------------------------------------------------

# -- "systems" (right aligned, wraps at box) --
$tbSystems = $newSlide.Shapes.AddTextbox(1, 14.2024, 192.9761, 359.8128, 94.5141)
$tbSystems.Name = "TextBox 1"
$tbSystems.TextFrame.WordWrap = -1
$tbSystems.TextFrame.AutoSize = 1
$tbSystems.TextFrame.TextRange.Text = "systems"
$tbSystems.TextFrame.TextRange.ParagraphFormat.Alignment = 3
$tbSystems.TextFrame.TextRange.Font.Name = "Georgia"
$tbSystems.TextFrame.TextRange.Font.Size = 72
$tbSystems.TextFrame.TextRange.Font.Bold = $true
$tbSystems.TextFrame.TextRange.Font.Color.RGB = HexColor("9CB1C7")

# -- "medicine" (right aligned, no-wrap autofit) --
$tbMedicine = $newSlide.Shapes.AddTextbox(1, 14.1348, 252.5097, 359.8805, 94.5141)
$tbMedicine.Name = "TextBox 2"
$tbMedicine.TextFrame.WordWrap = 0
$tbMedicine.TextFrame.AutoSize = 1
$tbMedicine.TextFrame.TextRange.Text = "medicine"
$tbMedicine.TextFrame.TextRange.ParagraphFormat.Alignment = 3
$tbMedicine.TextFrame.TextRange.Font.Name = "Georgia"
$tbMedicine.TextFrame.TextRange.Font.Size = 72
$tbMedicine.TextFrame.TextRange.Font.Bold = $true
$tbMedicine.TextFrame.TextRange.Font.Color.RGB = HexColor("9CB1C7")

# -- "disease" --
$tbDisease = $newSlide.Shapes.AddTextbox(1, 399.5323, 192.9761, 289.3231, 94.5141)
$tbDisease.Name = "TextBox 3"
$tbDisease.TextFrame.WordWrap = 0
$tbDisease.TextFrame.AutoSize = 1
$tbDisease.TextFrame.TextRange.Text = "disease"
$tbDisease.TextFrame.TextRange.Font.Name = "Georgia"
$tbDisease.TextFrame.TextRange.Font.Size = 72
$tbDisease.TextFrame.TextRange.Font.Bold = $true
$tbDisease.TextFrame.TextRange.Font.Color.RGB = HexColor("2E5287")

# -- "maps" --
$tbMaps = $newSlide.Shapes.AddTextbox(1, 399.5323, 252.5097, 216.7462, 94.5141)
$tbMaps.Name = "TextBox 4"
$tbMaps.TextFrame.WordWrap = 0
$tbMaps.TextFrame.AutoSize = 1
$tbMaps.TextFrame.TextRange.Text = "maps"
$tbMaps.TextFrame.TextRange.Font.Name = "Georgia"
$tbMaps.TextFrame.TextRange.Font.Size = 72
$tbMaps.TextFrame.TextRange.Font.Bold = $true
$tbMaps.TextFrame.TextRange.Font.Color.RGB = HexColor("2E5287")

# Group the four word-marks together ("Group 5").
$wordRange = $newSlide.Shapes.Range(@($tbSystems.Name, $tbMedicine.Name, $tbDisease.Name, $tbMaps.Name))
$group5 = $wordRange.Group()
$group5.Name = "Group 5"

# -- small rounded-rectangle accent next to the word-marks --
$roundRect = $newSlide.Shapes.AddShape(5, 384.6469, 213.3008, 5.6699, 124.7383)
$roundRect.Name = "Rounded Rectangle 6"
$roundRect.Adjustments.Item(1) = 0.47427
$roundRect.Fill.ForeColor.RGB = HexColor("9CB1C7")
$roundRect.Line.ForeColor.RGB = HexColor("9CB1C7")
$roundRect.TextFrame.TextRange.Text = ""

# Group the word-marks group with the rounded rectangle ("Group 7").
$outerRange = $newSlide.Shapes.Range(@($group5.Name, $roundRect.Name))
$group7 = $outerRange.Group()
$group7.Name = "Group 7"
